$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.365.09"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "1.843.96"
$ws.Range("E3").Value = "  -0.42%  "

Set-TextValue $ws.Range("D4") "0.9985"
$ws.Range("E4").Value = "  -0.08%  "

Set-TextValue $ws.Range("D5") "240.33"
$ws.Range("E5").Value = "  -1.13%  "

Set-TextValue $ws.Range("D6") "0.6346"
$ws.Range("E6").Value = "  -0.42%  "

Set-TextValue $ws.Range("D7") "0.9997"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "6.457.39"
$ws.Range("E8").Value = "  +248.74%  "

Set-TextValue $ws.Range("D9") "0.07544"
$ws.Range("E9").Value = "  -0.40%  "

Set-TextValue $ws.Range("D10") "0.2960"
$ws.Range("E10").Value = "  -1.55%  "

Set-TextValue $ws.Range("D11") "24.71"
$ws.Range("E11").Value = "  +2.03%  "

Set-TextValue $ws.Range("D12") "0.07733"
$ws.Range("E12").Value = "  +0.58%  "

Set-TextValue $ws.Range("D13") "4.988"
$ws.Range("E13").Value = "  -1.07%  "

Set-TextValue $ws.Range("D14") "0.6831"
$ws.Range("E14").Value = "  -0.78%  "

Set-TextValue $ws.Range("D15") "83.18"
$ws.Range("E15").Value = "  -1.08%  "

Set-TextValue $ws.Range("D16") "0.000009900"
$ws.Range("E16").Value = "  +1.52%  "

Set-TextValue $ws.Range("D17") "6.160"
$ws.Range("E17").Value = "  -2.30%  "

$ws.Range("D18").Value = "29.384.49"
$ws.Range("E18").Value = "  -0.82%  "

Set-TextValue $ws.Range("D19") "229.08"
$ws.Range("E19").Value = "  -4.09%  "

Set-TextValue $ws.Range("D20") "12.45"
$ws.Range("E20").Value = "  -0.82%  "

Set-TextValue $ws.Range("D21") "0.9996"
$ws.Range("E21").Value = "  -0.07%  "

Set-TextValue $ws.Range("D22") "7.545"
$ws.Range("E22").Value = "  -0.85%  "

Set-TextValue $ws.Range("D23") "1.000"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +234.54%  "

$ws.Range("E25").Value = "  +171.02%  "

Set-TextValue $ws.Range("D26") "155.94"
$ws.Range("E26").Value = "  -0.50%  "

Set-TextValue $ws.Range("D27") "0.1405"
$ws.Range("E27").Value = "  +0.59%  "

Set-TextValue $ws.Range("D28") "8.377"
$ws.Range("E28").Value = "  -0.91%  "

Set-TextValue $ws.Range("D29") "17.65"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("E30").Value = "  -1.08%  "

Set-TextValue $ws.Range("D31") "0.05705"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("E32").Value = "  -2.06%  "

Set-TextValue $ws.Range("D33") "4.125"
$ws.Range("E33").Value = "  -0.16%  "

Set-TextValue $ws.Range("D34") "4.030"
$ws.Range("E34").Value = "  -1.15%  "

Set-TextValue $ws.Range("D35") "1.844"
$ws.Range("E35").Value = "  -3.28%  "

Set-TextValue $ws.Range("D36") "1.156"
$ws.Range("E36").Value = "  -1.72%  "

Set-TextValue $ws.Range("D37") "0.7192"
$ws.Range("E37").Value = "  -0.33%  "

Set-TextValue $ws.Range("D38") "2.590"
$ws.Range("E38").Value = "  -0.33%  "

$ws.Range("D39").Value = "1.250.99"
$ws.Range("E39").Value = "  +1.83%  "

Set-TextValue $ws.Range("D40") "2.800"
$ws.Range("E40").Value = "  -0.32%  "

Set-TextValue $ws.Range("D41") "0.01812"
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.9065"
$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "0.9994"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D44") "102.04"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "66.37"
$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "7.063"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "9.166"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D48") "0.4022"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "1.706"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.1125"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05738"
$ws.Range("E51").Value = "  -0.25%  "
